# This weekly data extract gained one new observation row. The source feed
# prepends newer records, which pushes every existing record for this
# subset down by one row. Reproduce that by inserting a fresh row 20 (Excel
# shifts rows 20:116 down to 21:117 automatically, carrying their values and
# formatting with them) and then filling the new row with the incoming
# weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20; everything currently at/after row 20
# (through the former last row 116) shifts down to 21..117.
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new weekly record.
$ws.Range("A20").Value = 7
$ws.Range("B20").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C20").Value = "Ñuble"
$ws.Range("D20").Value = 44547
$ws.Range("E20").Value = 16
$ws.Range("F20").Value = 100112024
$ws.Range("G20").Value = "Choclo"
$ws.Range("H20").Value = "Choclero"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 12000
$ws.Range("K20").Value = 350
$ws.Range("L20").Value = 400
$ws.Range("M20").Value = 375
$ws.Range("N20").Value = "$/unidad"
$ws.Range("O20").Value = "Región de O'Higgins"
$ws.Range("P20").Value = 375
$ws.Range("Q20").Value = 1
$ws.Range("R20").Value = "Hortaliza"
